# Add a new "AutoAtkDis" property/column for skills.
# This inserts a new column before the existing "NeedTar" column (column I),
# shifting AtkDis's neighboring columns (NeedTar, DefaultHitTime, ShowName)
# one slot to the right, then populates the header and fills the new column
# with the default value 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (pushes NeedTar/DefaultHitTime/ShowName right).
$ws.Columns("I:I").Insert()

# New column header + data.
$ws.Range("I1").Value = "AutoAtkDis"
$ws.Range("I2:I9").Value = 1

# Leave the selection where the edit happened, like Excel would after
# filling the new column in by hand.
$ws.Range("I2:I9").Select() | Out-Null
